$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Course")

# New course row (row 4): CourseID, Course Name, Course Link, Course Author, Certificate
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "The Complete Financial Analyst Course 2023"
$ws.Range("C4").Value = "https://www.udemy.com/course/the-complete-financial-analyst-course"
$ws.Range("D4").Value = "365 Careers"
$ws.Range("E4").Value = "No"

# Hyperlink on the course-link cell (added before style copy so the copy wins on formatting)
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.udemy.com/course/the-complete-financial-analyst-course")

# Reapply formatting from the row above so new cells share the same styles
# (Hyperlink style for C4, date format for F4/J4) instead of Excel
# auto-generating brand-new style records.
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = 45057

$ws.Range("J3").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 45057

$excel.CutCopyMode = $false

# Match the saved selection state
$ws.Range("J5").Select()
